$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1): make bold (content stays the same)
# ---------------------------------------------------------------------------
$ws.Range("A1:E1").Font.Bold = $true

# ---------------------------------------------------------------------------
# 2. Remove the old "test" intent rows (rows 2,4,5,6,7,8) and replace the
#    first block (rows 2-3) with the new "buy / product.category" intents.
#    Unmerge the groups that are being dissolved.
# ---------------------------------------------------------------------------
$ws.Range("A2:A3").UnMerge()
$ws.Range("B2:B3").UnMerge()
$ws.Range("C2:C3").UnMerge()
$ws.Range("D2:D3").UnMerge()
$ws.Range("A4:A5").UnMerge()
$ws.Range("B4:B5").UnMerge()
$ws.Range("A6:A8").UnMerge()
$ws.Range("B6:B8").UnMerge()

# Clear the whole block (rows 2-8, columns A-E) before rewriting it
$ws.Range("A2:E8").ClearContents()

# New row 2 : "buy>product.category" intent
$ws.Range("A2").Value = "buy>product.category"
$ws.Range("B2").Value = "`ncomprar`nquiero hacer una compra`nestoy buscando`nquiero comprar"

# New row 3 : "product.category>sp.range" intent
$ws.Range("A3").Value = "product.category>sp.range"
$ws.Range("B3").Value = "quiero comprar un smartphone`nestoy buscando un m$([char]0xF3)vil`nno s$([char]0xE9) que m$([char]0xF3)vil comprarme`nme gustar$([char]0xED)a comprar un smartphone`nayudame a elegir un movil"
$ws.Range("C3").Value = "productCategory"
$ws.Range("D3").Value = "smartphone`nm$([char]0xF3)vil`nmoviles`nmovil"

# ---------------------------------------------------------------------------
# 3. Formatting / alignment for rows 2-14, columns A-E : centered horizontally
#    and vertically; B/D columns (the "UserSays"/"AnnotationValue" columns)
#    additionally wrap text.
# ---------------------------------------------------------------------------
$ws.Range("A2:E14").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A2:E14").VerticalAlignment = -4108     # xlCenter

$ws.Range("B2:B3").WrapText = $true
$ws.Range("D2:D3").WrapText = $true
$ws.Range("B9").WrapText = $true
$ws.Range("B11").WrapText = $true
$ws.Range("B12:B14").WrapText = $true
$ws.Range("D12:D14").WrapText = $true

# ---------------------------------------------------------------------------
# 4. Row heights / dimensions for the new layout
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 76.5
$ws.Rows.Item(3).RowHeight = 75
$ws.Rows.Item(8).RowHeight = 76.5

# ---------------------------------------------------------------------------
# 5. Column widths (A and E change size)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 32
$ws.Columns.Item(5).ColumnWidth = 25.65

# ---------------------------------------------------------------------------
# 6. Row 15 gains vertical-centered (no horizontal) formatting across A:D
# ---------------------------------------------------------------------------
$ws.Range("A15:D15").HorizontalAlignment = -4142   # xlGeneral (no horizontal centering)
$ws.Range("A15:D15").VerticalAlignment = -4108     # xlCenter

# ---------------------------------------------------------------------------
# 7. Selection moves to D3
# ---------------------------------------------------------------------------
$ws.Range("D3").Select()
